# Fixes on single board BOM sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keyboard backlight row: add JLCPCB part number.
$ws.Range("D4").Value = "C2895458"

# 1206 Orange LED row: JLCPCB part number corrected.
$ws.Range("D8").Value = "C2764895"

# Keyboard switches (MX) row: add JLCPCB part number.
$ws.Range("D9").Value = "C5120587"

# Reverse polarity protection: footprint/diode changed, JLCPCB part number corrected.
$ws.Range("C13").Value = "MBR0520LT1G"
$ws.Range("D13").Value = "C23848"
$ws.Range("C13").WrapText = $true

# Update the active selection to D2.
$ws.Range("D2").Select() | Out-Null
